$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.959.21'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.617.09'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '211.65'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.488'
$ws.Range("E7").Value = '  -3.35%  '
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '18.16'
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D12").Value = '1.841.02'
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '4.12'
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.594.80'
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.518'
$ws.Range("E15").Value = '  -2.24%  '
$ws.Range("D16").Value = '25.960.98'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '61.53'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = '0.0₃0731'
$ws.Range("E18").Value = '  -2.08%  '
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '191.17'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '4.22'
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '9.41'
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '6.00'
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '0.131'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '143.10'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '1.73'
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '6.60'
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '15.13'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '1.22'
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  -2.35%  '
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '3.07'
$ws.Range("E33").Value = '  -2.69%  '
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("D36").Value = '1.121.53'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.818'
$ws.Range("E37").Value = '  -6.52%  '
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("E39").Value = '  -2.41%  '
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '97.17'
$ws.Range("E41").Value = '  -1.64%  '
$ws.Range("D42").Value = '1.751.62'
$ws.Range("E42").Value = '  -1.11%  '
$ws.Range("E43").Value = '  -3.96%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '5.05'
$ws.Range("E44").Value = '  -4.82%  '
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '53.74'
$ws.Range("E46").Value = '  -2.93%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '1.48'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.410'
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '7.41'
$ws.Range("E51").Value = '  -1.47%  '

Write-Output "Done applying cryptos update"
